# Insert a new weekly price record for Perejil (Terminal La Palmera de La
# Serena) at row 98. This pushes the existing rows 98-130 down to 99-131
# (the previously-last row, dated 44544, ends up at row 131), and fills the
# freshly inserted row 98 with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(98).Insert()

$ws.Range("A98").Value = 8
$ws.Range("B98").Value = "Terminal La Palmera de La Serena"
$ws.Range("C98").Value = "Coquimbo"
$ws.Range("D98").Value = 44642
$ws.Range("E98").Value = 4
$ws.Range("F98").Value = 100112044
$ws.Range("G98").Value = "Perejil"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 2400
$ws.Range("K98").Value = 2300
$ws.Range("L98").Value = 2500
$ws.Range("M98").Value = 2400
$ws.Range("N98").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O98").Value = "Provincia del Elquí"
$ws.Range("P98").Value = 1600
$ws.Range("Q98").Value = 1.5
$ws.Range("R98").Value = "Hortaliza"
